$wb = $excel.ActiveWorkbook

# Use a throwaway sheet to burn a sheetId so the new sheet lands on sheetId=10
# (Excel's id counter is "max currently-used id + 1" at the time a sheet is
# created, so we create-and-discard one placeholder first).
$placeholder = $wb.Worksheets.Add()
$placeholder.Name = "__placeholder__"

# The target sheet we are duplicating/inserting before. Fetched by name
# *after* the placeholder insert so it resolves to the right worksheet.
$src = $wb.Worksheets.Item("Allocation_Summary")

# Copy Allocation_Summary (with all formatting/content) to just before itself.
$src.Copy($src)

# Remove the placeholder now that the real copy owns sheetId=10.
$wb.Worksheets.Item("__placeholder__").Delete()

# The copy was named "Allocation_Summary (2)"; rename + re-point its content.
$newSheet = $wb.Worksheets.Item("Allocation_Summary (2)")
$newSheet.Name = "CA-AllocationSummary"
$newSheet.Range("A2").Value = "CA-AllocationSummary"

# Make the new sheet the active/selected tab with A7 as the active cell.
$newSheet.Select()
$newSheet.Range("A7").Select()
